$d = $word.ActiveDocument

$pairs = @(
    @("496×4=", "571×9="),
    @("736×4=", "355×7="),
    @("573×3=", "693×7="),
    @("163×8=", "829×7="),
    @("748×3=", "726×7="),
    @("434×3=", "551×4="),
    @("859×6=", "394×6="),
    @("644×8=", "367×5="),
    @("677×2=", "562×3="),
    @("784×6=", "899×2="),
    @("408×3=", "241×7="),
    @("342×2=", "914×2="),
    @("516×2=", "161×6="),
    @("233×9=", "642×3="),
    @("833×4=", "714×7="),
    @("521×5=", "560×7="),
    @("419×4=", "231×3="),
    @("237×8=", "898×9="),
    @("561×9=", "342×5="),
    @("185×6=", "574×6="),
    @("506×5=", "991×3="),
    @("632×6=", "969×8="),
    @("145×3=", "980×9="),
    @("260×3=", "188×3="),
    @("144×8=", "602×6=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
